$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date style used by column D in existing data rows (style index "2" / custom date number format)
$dateFormat = $ws.Cells.Item(2, 4).NumberFormat

$rows = @(
    @{ Row = 131; K = "Fuji royal";    L = "Calibre 80" },
    @{ Row = 132; K = "Granny Smith";  L = "Calibre 90" },
    @{ Row = 133; K = "Royal Gala";    L = "Calibre 90" }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value2 = 1
    $ws.Cells.Item($row, 2).Value2 = "Agrícola del Norte S.A. de Arica"
    $ws.Cells.Item($row, 3).Value2 = "Arica y Parinacota"

    $ws.Cells.Item($row, 4).Value2 = 44747
    $ws.Cells.Item($row, 4).NumberFormat = $dateFormat

    $ws.Cells.Item($row, 5).Value2 = 15
    $ws.Cells.Item($row, 6).Value2 = "Fruta"
    $ws.Cells.Item($row, 7).Value2 = 100104
    $ws.Cells.Item($row, 8).Value2 = "Frutos de pepita"
    $ws.Cells.Item($row, 9).Value2 = 100104002
    $ws.Cells.Item($row, 10).Value2 = "Manzana"
    $ws.Cells.Item($row, 11).Value2 = $r.K
    $ws.Cells.Item($row, 12).Value2 = $r.L
    $ws.Cells.Item($row, 13).Value2 = 300
    $ws.Cells.Item($row, 14).Value2 = 17000
    $ws.Cells.Item($row, 15).Value2 = 18000
    $ws.Cells.Item($row, 16).Value2 = 17500
    $ws.Cells.Item($row, 17).Value2 = "$/caja 18 kilos embalada"
    $ws.Cells.Item($row, 18).Value2 = "Región de O'Higgins"
    $ws.Cells.Item($row, 19).Value2 = 972
    $ws.Cells.Item($row, 20).Value2 = 18
}
